# Insert a new data row at row 17 (pushing all existing data rows 17..135
# down to 18..136), then populate the new row 17 with a new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 17; this shifts rows 17-135 down to 18-136
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new record values.
$ws.Cells.Item(17, 1).Value2 = 4
$ws.Cells.Item(17, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(17, 3).Value2 = "Los Lagos"
$ws.Cells.Item(17, 4).Value2 = 44473
$ws.Cells.Item(17, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value2 = 10
$ws.Cells.Item(17, 6).Value2 = 100112017
$ws.Cells.Item(17, 7).Value2 = "Apio"
$ws.Cells.Item(17, 8).Value2 = "Americana (o)"
$ws.Cells.Item(17, 9).Value2 = "Primera"
$ws.Cells.Item(17, 10).Value2 = 35
$ws.Cells.Item(17, 11).Value2 = 12000
$ws.Cells.Item(17, 12).Value2 = 12000
$ws.Cells.Item(17, 13).Value2 = 12000
$ws.Cells.Item(17, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(17, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(17, 16).Value2 = 2000
$ws.Cells.Item(17, 17).Value2 = 6
$ws.Cells.Item(17, 18).Value2 = "Hortaliza"
